$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-484) holds the "Förändrad" date, stored as Excel serial date 45171
# (2023-09-02). Every row's value is bumped by one day to 45172 (2023-09-03).
$lastRow = 484
$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45172
